$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.740.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.48%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.534.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.59%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'309.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.17%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'102.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.52%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -1.16%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.12%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.42%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'36.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.15%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.38%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -2.44%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.08%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.929.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.45%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.69%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.488.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.27%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.812"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -4.46%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.737.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.83%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -2.04%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.77%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'12.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.61%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'69.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.49%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'245.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.98%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -2.66%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -1.86%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.03%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'26.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -5.11%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -3.99%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'39.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.63%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.31%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'157.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.29%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -2.00%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'2.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +9.25%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -2.37%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.78%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -6.04%  "
$ws.Range("E36").Style = "Normal"
# Rows 37 and 38: coins swapped (Celestia <-> LidoDAOToken) with updated price/volume
$ws.Range("B37").Value = "'LidoDAOToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'3.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -9.45%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'Celestia"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'18.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.59%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -0.44%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.06%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'4.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +6.79%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'22.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.13%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.00%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'3.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.62%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.41%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.985.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.82%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.29%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.784.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.53%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'80.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.80%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.09%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.853"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +8.15%  "
$ws.Range("E51").Style = "Normal"
